# Generate Report for Handback
# Update the Correspond Handoff/Handback datetime stamps for the
# 60707580-db3e-4393-861b-1229ed4c4f94 entry (row 2 / row 4, same values)
# on both the "zh-cn" and "de-de" language sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 06:18:28"
$wsZhCn.Range("H2").Value = "2016-03-19 06:18:48"
$wsZhCn.Range("E4").Value = "2016-03-19 06:18:28"
$wsZhCn.Range("H4").Value = "2016-03-19 06:18:48"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 06:18:32"
$wsDeDe.Range("H2").Value = "2016-03-19 06:18:53"
$wsDeDe.Range("E4").Value = "2016-03-19 06:18:32"
$wsDeDe.Range("H4").Value = "2016-03-19 06:18:53"
